$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit effectively swaps the full contents of row 2 with row 4,
# and the full contents of row 3 with row 5 (a cyclic re-ordering of
# the 4 observation records). Only swap the columns that actually hold
# data for these rows, to avoid disturbing already-blank cells.

function Swap-Cell($ws, $addr1, $addr2) {
    $v1 = $ws.Range($addr1).Value2
    $v2 = $ws.Range($addr2).Value2
    $ws.Range($addr1).Value2 = $v2
    $ws.Range($addr2).Value2 = $v1
}

# Columns with regular (non date-like) values that differ between row 2 and row 4.
$cols24 = @("A","B","D","E","F","G","H","Q","R","S")
foreach ($col in $cols24) {
    $addr1 = $col + "2"
    $addr2 = $col + "4"
    Swap-Cell $ws $addr1 $addr2
}

# Columns with regular values that differ between row 3 and row 5.
$cols35 = @("A","B","D","E","F","G","H","Q","R","Z","AB")
foreach ($col in $cols35) {
    $addr1 = $col + "3"
    $addr2 = $col + "5"
    Swap-Cell $ws $addr1 $addr2
}

# Column I holds a text value that looks numeric ("35"). A plain Value2
# assignment would be auto-converted to a real number, so instead stage
# the text in a scratch cell formatted as Text, copy/paste it as a value
# (this preserves the text type), and finally wipe the scratch cell so
# it leaves no trace in the saved workbook.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value2 = "35"
$scratch.Copy()
$ws.Range("I3").PasteSpecial(-4163) # xlPasteValues
$ws.Range("I5").ClearContents()
$scratch.Clear()

$excel.CutCopyMode = 0
